# Auto-generated Excel COM-interop script
# Applies a scheduled market-data refresh to the price/profit columns (H:N)
# across the 8 crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values come from an external market API snapshot; columns are plain numbers
# (no formulas in this workbook), so each target cell is set/cleared directly.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2394.8262
$ws.Range("I40").Value = 2260.6924
$ws.Range("K40").Value = 2260.6924
$ws.Range("M40").Value = -2085.6924

$ws.Range("H53").Value = 317.6
$ws.Range("I53").Value = 265.73334
$ws.Range("K53").Value = 265.73334
$ws.Range("M53").Value = 371.26666

$ws.Range("H74").Value = 49476.895
$ws.Range("I74").Value = 70184.836
$ws.Range("J74").Value = 12202.6
$ws.Range("K74").Value = 70184.836
$ws.Range("L74").Value = 12202.6
$ws.Range("M74").Value = -69248.836
$ws.Range("N74").Value = -14074.6

$ws.Range("H77").Value = 49476.895
$ws.Range("I77").Value = 70184.836
$ws.Range("J77").Value = 12202.6
$ws.Range("K77").Value = 350924.18
$ws.Range("L77").Value = 61013
$ws.Range("M77").Value = -346244.18
$ws.Range("N77").Value = -70373

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3215.25
$ws.Range("I31").Value = 3215.25
$ws.Range("K31").Value = 3215.25
$ws.Range("M31").Value = -2921.25

$ws.Range("H61").Value = 6262.091
$ws.Range("I61").Value = 1630.3334
$ws.Range("K61").Value = 1630.3334
$ws.Range("M61").Value = -1418.3334

$ws.Range("H132").Value = 1454.1
$ws.Range("I132").Value = 1321.7059
$ws.Range("J132").Value = 2204.3333
$ws.Range("K132").Value = 3965.1177
$ws.Range("L132").Value = 6612.999899999999
$ws.Range("M132").Value = -1435.1177
$ws.Range("N132").Value = -11672.9999

$ws.Range("H136").Value = 6262.091
$ws.Range("I136").Value = 1630.3334
$ws.Range("K136").Value = 4891.0002
$ws.Range("M136").Value = -2341.0002


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3363.5833
$ws.Range("I99").Value = 2231.625
$ws.Range("K99").Value = 2231.625
$ws.Range("M99").Value = -733.625

$ws.Range("H107").Value = 3613.4167
$ws.Range("I107").Value = 3213.4546
$ws.Range("K107").Value = 3213.4546
$ws.Range("M107").Value = -1293.4546

$ws.Range("H115").Value = 48996.5
$ws.Range("J115").Value = 48996.5
$ws.Range("L115").Value = 48996.5
$ws.Range("N115").Value = -52130.5


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1577
$ws.Range("I16").Value = 1577
$ws.Range("K16").Value = 1577
$ws.Range("M16").Value = -1290

$ws.Range("H29").Value = 1274.75
$ws.Range("J29").Value = 1366.6666
$ws.Range("L29").Value = 1366.6666
$ws.Range("N29").Value = -1952.6666

$ws.Range("H31").Value = 2096.36
$ws.Range("I31").Value = 1755.619
$ws.Range("K31").Value = 1755.619
$ws.Range("M31").Value = -1460.619

$ws.Range("H34").Value = 2096.36
$ws.Range("I34").Value = 1755.619
$ws.Range("K34").Value = 1755.619
$ws.Range("M34").Value = -1553.619

$ws.Range("H51").Value = 30045
$ws.Range("I51").Value = 30090
$ws.Range("J51").Value = 30000
$ws.Range("K51").Value = 30090
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = -29354
$ws.Range("N51").Value = -31472

$ws.Range("H61").Value = 30045
$ws.Range("I61").Value = 30090
$ws.Range("J61").Value = 30000
$ws.Range("K61").Value = 30090
$ws.Range("L61").Value = 30000
$ws.Range("M61").Value = -29742
$ws.Range("N61").Value = -30696

$ws.Range("H113").Value = 1577
$ws.Range("I113").Value = 1577
$ws.Range("K113").Value = 1577
$ws.Range("M113").Value = 593

$ws.Range("H122").Value = 1289.4166
$ws.Range("I122").Value = 1041.7273
$ws.Range("K122").Value = 3125.1819
$ws.Range("M122").Value = -675.1819


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1344.3334
$ws.Range("I14").Value = 1344.3334
$ws.Range("K14").Value = 4033.0002
$ws.Range("M14").Value = -3860.0002

$ws.Range("H25").Value = 800500.2
$ws.Range("I25").Value = 2501
$ws.Range("J25").Value = 1000000
$ws.Range("K25").Value = 7503
$ws.Range("L25").Value = 3000000
$ws.Range("M25").Value = -7334
$ws.Range("N25").Value = -3000338

$ws.Range("H30").Value = 800500.2
$ws.Range("I30").Value = 2501
$ws.Range("J30").Value = 1000000
$ws.Range("K30").Value = 7503
$ws.Range("L30").Value = 3000000
$ws.Range("M30").Value = -7401
$ws.Range("N30").Value = -3000204

$ws.Range("H35").Value = 325
$ws.Range("I35").Value = 325
$ws.Range("K35").Value = 975
$ws.Range("M35").Value = -687

$ws.Range("H49").Value = 2000
$ws.Range("J49").Value = 2000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6312

$ws.Range("H58").Value = 2495
$ws.Range("J58").Value = 2495
$ws.Range("L58").Value = 7485
$ws.Range("N58").Value = -7741

$ws.Range("H59").Value = 495
$ws.Range("J59").Value = 495
$ws.Range("L59").Value = 1485
$ws.Range("N59").Value = -2565

$ws.Range("H68").Value = 20847332
$ws.Range("I68").Value = 20999
$ws.Range("K68").Value = 62997
$ws.Range("M68").Value = -62186

$ws.Range("H71").Value = 20847332
$ws.Range("I71").Value = 20999
$ws.Range("K71").Value = 188991
$ws.Range("M71").Value = -184935

$ws.Range("H110").Value = 5027
$ws.Range("I110").Value = 5027
$ws.Range("K110").Value = 15081
$ws.Range("M110").Value = -10991

$ws.Range("H139").Value = 4705.4
$ws.Range("I139").Value = 4705.4
$ws.Range("K139").Value = 14116.2
$ws.Range("M139").Value = -8976.199999999999


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 875000.75
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H14").Value = 6001417.5
$ws.Range("I14").Value = 7201100
$ws.Range("J14").Value = 3005
$ws.Range("K14").Value = 7201100
$ws.Range("L14").Value = 3005
$ws.Range("M14").Value = -7200932
$ws.Range("N14").Value = -3341

$ws.Range("H80").Value = 3103.8667
$ws.Range("J80").Value = 3119.4614
$ws.Range("L80").Value = 3119.4614
$ws.Range("N80").Value = -5115.4614

$ws.Range("H83").Value = 3103.8667
$ws.Range("J83").Value = 3119.4614
$ws.Range("L83").Value = 15597.307
$ws.Range("N83").Value = -25581.307


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3725.0715
$ws.Range("I40").Value = 2621.7144
$ws.Range("K40").Value = 2621.7144
$ws.Range("M40").Value = -2485.7144

$ws.Range("H61").Value = 4151.8696
$ws.Range("I61").Value = 3972.7778
$ws.Range("K61").Value = 3972.7778
$ws.Range("M61").Value = -3770.7778

$ws.Range("H82").Value = 719.8
$ws.Range("I82").Value = 714.125
$ws.Range("J82").Value = 742.5
$ws.Range("K82").Value = 714.125
$ws.Range("L82").Value = 742.5
$ws.Range("M82").Value = -353.125
$ws.Range("N82").Value = -1464.5

$ws.Range("H85").Value = 719.8
$ws.Range("I85").Value = 714.125
$ws.Range("J85").Value = 742.5
$ws.Range("K85").Value = 714.125
$ws.Range("L85").Value = 742.5
$ws.Range("M85").Value = 533.875
$ws.Range("N85").Value = -3238.5

$ws.Range("H113").Value = 4151.8696
$ws.Range("I113").Value = 3972.7778
$ws.Range("K113").Value = 3972.7778
$ws.Range("M113").Value = -1802.7778

$ws.Range("H132").Value = 1784.2858
$ws.Range("I132").Value = 1815
$ws.Range("K132").Value = 5445
$ws.Range("M132").Value = -2915


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 599
$ws.Range("I107").Value = 599
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1797
$ws.Range("L107").ClearContents()
$ws.Range("N107").Value = 0
$ws.Range("M107").Value = 123

$ws.Range("H113").Value = 548.2222
$ws.Range("I113").Value = 187
$ws.Range("J113").Value = 999.75
$ws.Range("K113").Value = 561
$ws.Range("L113").Value = 2999.25
$ws.Range("M113").Value = 1609
$ws.Range("N113").Value = -7339.25

$ws.Range("H132").Value = 1595.1818
$ws.Range("I132").Value = 1355.2
$ws.Range("K132").Value = 4065.6
$ws.Range("M132").Value = -1535.6

